$wb = $excel.ActiveWorkbook

# --- 1. Duplicate the "demands" sheet to create "contents" -----------------
# (the old "demands" data - id/demand_model_id/resource_id/amount/rate - is
# kept verbatim on the new "contents" tab, placed just before "demand_models")
$demandsOriginal = $wb.Worksheets.Item("demands")
$demandModels = $wb.Worksheets.Item("demand_models")
$demandsOriginal.Copy($demandModels)

$contents = $wb.Worksheets.Item("demands (2)")
$contents.Name = "contents"
$contents = $wb.Worksheets.Item("contents")
$contents.Activate()
$contents.Columns("D").Select()

# --- 2. Replace the content of the (original) "demands" sheet --------------
# with the new container-based schema: id, container_id, resource_id, amount
$demands = $wb.Worksheets.Item("demands")
$demands.Cells.Clear()
$demands.Range("A1").Value = "id"
$demands.Range("B1").Value = "container_id"
$demands.Range("C1").Value = "resource_id"
$demands.Range("D1").Value = "amount"
$demands.Activate()
$demands.Range("G20").Select()

# --- 3. Misc selection / active-sheet bookkeeping to match the edit --------
$demandModels = $wb.Worksheets.Item("demand_models")
$demandModels.Activate()
$demandModels.Range("G48").Select()

$resources = $wb.Worksheets.Item("resources")
$resources.Activate()
$resources.Range("I7").Select()

$states = $wb.Worksheets.Item("states")
$states.Activate()
$states.Range("C18").Select()

$elements = $wb.Worksheets.Item("elements")
$elements.Activate()
$elements.Range("B8").Select()
